$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure columns D and E keep their text formatting so that numeric-looking
# strings (e.g. "1.20", "0.0000247", "98.137.42") are not coerced into numbers
# by Excel when the new values are assigned.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "98.137.42"
$ws.Range("E2").Value = "  +4.07%  "

$ws.Range("D3").Value = "3.369.48"
$ws.Range("E3").Value = "  +9.20%  "

$ws.Range("E4").Value = "  -0.06%  "

$ws.Range("D5").Value = "254.33"
$ws.Range("E5").Value = "  +7.29%  "

$ws.Range("D6").Value = "622.67"
$ws.Range("E6").Value = "  +1.98%  "

$ws.Range("D7").Value = "1.20"
$ws.Range("E7").Value = "  +7.87%  "

$ws.Range("D8").Value = "0.385"
$ws.Range("E8").Value = "  +1.44%  "

$ws.Range("E9").Value = "  -0.02%  "

$ws.Range("D10").Value = "3.364.59"
$ws.Range("E10").Value = "  +9.17%  "

$ws.Range("D11").Value = "0.806"
$ws.Range("E11").Value = "  -0.86%  "

$ws.Range("E12").Value = "  +1.18%  "

$ws.Range("D13").Value = "97.817.91"
$ws.Range("E13").Value = "  +4.00%  "

$ws.Range("D14").Value = "35.84"
$ws.Range("E14").Value = "  +4.91%  "

$ws.Range("D15").Value = "0.0000247"
$ws.Range("E15").Value = "  +2.16%  "

$ws.Range("D16").Value = "3.991.69"
$ws.Range("E16").Value = "  +9.03%  "

$ws.Range("D17").Value = "5.49"
$ws.Range("E17").Value = "  +2.85%  "

$ws.Range("D18").Value = "3.369.72"
$ws.Range("E18").Value = "  +9.70%  "

$ws.Range("D19").Value = "3.65"
$ws.Range("E19").Value = "  +1.98%  "

$ws.Range("D20").Value = "14.79"
$ws.Range("E20").Value = "  +2.12%  "

$ws.Range("D21").Value = "480.86"
$ws.Range("E21").Value = "  +7.24%  "

$ws.Range("D22").Value = "5.89"
$ws.Range("E22").Value = "  +1.86%  "

$ws.Range("D23").Value = "0.0000208"
$ws.Range("E23").Value = "  +9.27%  "

$ws.Range("D24").Value = "9.17"
$ws.Range("E24").Value = "  +3.43%  "

$ws.Range("D25").Value = "5.71"
$ws.Range("E25").Value = "  +3.12%  "

$ws.Range("D26").Value = "88.09"
$ws.Range("E26").Value = "  +3.73%  "

$ws.Range("D27").Value = "12.03"
$ws.Range("E27").Value = "  +0.05%  "

$ws.Range("D28").Value = "3.550.85"
$ws.Range("E28").Value = "  +9.45%  "

$ws.Range("E29").Value = "  -0.16%  "

$ws.Range("D30").Value = "0.188"
$ws.Range("E30").Value = "  +4.93%  "

$ws.Range("E31").Value = "  -0.35%  "

$ws.Range("D32").Value = "0.126"
$ws.Range("E32").Value = "  +2.45%  "

$ws.Range("D33").Value = "0.999"
$ws.Range("E33").Value = "  +0.13%  "

$ws.Range("D34").Value = "9.26"
$ws.Range("E34").Value = "  +2.53%  "

$ws.Range("D35").Value = "27.38"
$ws.Range("E35").Value = "  +7.16%  "

$ws.Range("D36").Value = "525.23"
$ws.Range("E36").Value = "  +8.98%  "

$ws.Range("D37").Value = "0.152"
$ws.Range("E37").Value = "  +0.25%  "

$ws.Range("D38").Value = "7.31"
$ws.Range("E38").Value = "  -4.62%  "

$ws.Range("D39").Value = "1.94"
$ws.Range("E39").Value = "  +2.85%  "

$ws.Range("D40").Value = "24.82"
$ws.Range("E40").Value = "  +3.17%  "

$ws.Range("D41").Value = "0.449"
$ws.Range("E41").Value = "  +2.23%  "

$ws.Range("B42").Value = "MantraDAO"
$ws.Range("C42").Value = "https://coinranking.com/coin/cTdD8lD-6+mantradao-om"
$ws.Range("D42").Value = "3.81"
$ws.Range("E42").Value = "  +2.56%  "

$ws.Range("B43").Value = "Fetch.AI"
$ws.Range("C43").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D43").Value = "1.27"
$ws.Range("E43").Value = "  +1.26%  "

$ws.Range("B44").Value = "dogwifhat"
$ws.Range("C44").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D44").Value = "3.25"
$ws.Range("E44").Value = "  +4.86%  "

$ws.Range("B45").Value = "ARBITRUM"
$ws.Range("C45").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D45").Value = "0.789"
$ws.Range("E45").Value = "  +16.74%  "

$ws.Range("E46").Value = "  -0.02%  "

$ws.Range("D47").Value = "161.18"
$ws.Range("E47").Value = "  -0.11%  "

$ws.Range("D48").Value = "1.93"
$ws.Range("E48").Value = "  +5.67%  "

$ws.Range("D49").Value = "45.52"
$ws.Range("E49").Value = "  +4.23%  "

$ws.Range("D50").Value = "1.37"
$ws.Range("E50").Value = "  +5.90%  "

$ws.Range("D51").Value = "4.55"
$ws.Range("E51").Value = "  +6.42%  "
